# Update handback status timestamps (report regeneration) as described by
# the commit "Generate Report for Handback".
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-09-06 11:26:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-09-06 11:26:11"
$wsZhCn.Range("K2").Value = "2016-09-06 11:26:33"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first file row.
$wsDeDe.Range("H2").Value = "2016-09-06 11:26:19"
$wsDeDe.Range("K2").Value = "2016-09-06 11:26:40"
